# Auto-generated: applies odds updates from the 2025-04-28 FlashScore refresh.
# Each assignment sets a cell to its new value per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.45
$ws.Range("H4").Value = 2.75
$ws.Range("I5").Value = 3.9
$ws.Range("L5").Value = 1.53
$ws.Range("M5").Value = 2.38
$ws.Range("P5").Value = 1.57
$ws.Range("Q5").Value = 2.25
$ws.Range("U5").Value = 8
$ws.Range("AA5").Value = 6.5
$ws.Range("AE5").Value = 8
$ws.Range("L6").Value = 1.44
$ws.Range("M6").Value = 2.63
$ws.Range("P9").Value = 1.53
$ws.Range("Q9").Value = 2.38
$ws.Range("P10").Value = 1.57
$ws.Range("P11").Value = 1.67
$ws.Range("G14").Value = 1.85
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 4.4
$ws.Range("K14").Value = 5.6
$ws.Range("N14").Value = 2.5
$ws.Range("O14").Value = 1.47
$ws.Range("Q14").Value = 2.32
$ws.Range("T14").Value = 5.2
$ws.Range("U14").Value = 7.3
$ws.Range("W14").Value = 14.5
$ws.Range("X14").Value = 18.5
$ws.Range("Z14").Value = 5.6
$ws.Range("AA14").Value = 6.4
$ws.Range("AB14").Value = 21
$ws.Range("AE14").Value = 9
$ws.Range("AF14").Value = 22
$ws.Range("AG14").Value = 16
$ws.Range("AI14").Value = 55
$ws.Range("I24").Value = 27
$ws.Range("N24").Value = 1.47
$ws.Range("O24").Value = 2.32
$ws.Range("R24").Value = 2.45
$ws.Range("T24").Value = 7.2
$ws.Range("X24").Value = 11.5
$ws.Range("Z24").Value = 13.5
$ws.Range("AE24").Value = 70
$ws.Range("AF24").Value = 400
$ws.Range("AG24").Value = 100
$ws.Range("AI24").Value = 800
$ws.Range("AJ24").Value = 400
$ws.Range("J30").Value = 1.05
$ws.Range("K30").Value = 11
$ws.Range("L30").Value = 1.29
$ws.Range("M30").Value = 3.5
$ws.Range("N30").Value = 1.93
$ws.Range("O30").Value = 1.88
$ws.Range("G31").Value = 2.85
$ws.Range("H31").Value = 2.45
$ws.Range("I31").Value = 3.15
$ws.Range("W31").Value = 37
$ws.Range("G33").Value = 3.35
$ws.Range("I33").Value = 2.67
$ws.Range("R33").Value = 2.12
$ws.Range("S33").Value = 1.65
$ws.Range("U33").Value = 16
$ws.Range("V33").Value = 12
$ws.Range("AF33").Value = 12
$ws.Range("AH33").Value = 35
$ws.Range("G34").Value = 1.53
$ws.Range("H34").Value = 4.1
$ws.Range("I34").Value = 6
$ws.Range("L34").Value = 1.22
$ws.Range("M34").Value = 4.33
$ws.Range("N34").Value = 1.73
$ws.Range("O34").Value = 2.1
$ws.Range("AA34").Value = 8
$ws.Range("AG34").Value = 19
$ws.Range("J39").Value = 1.08
$ws.Range("K39").Value = 8
$ws.Range("R39").Value = 2.63
$ws.Range("S39").Value = 1.44
$ws.Range("T39").Value = 5
$ws.Range("V39").Value = 9.5
$ws.Range("Z39").Value = 8
$ws.Range("AA39").Value = 9.5
$ws.Range("AC39").Value = 126
$ws.Range("AG39").Value = 29
$ws.Range("AI39").Value = 81
$ws.Range("AJ39").Value = 81
$ws.Range("G42").Value = 2.45
$ws.Range("I42").Value = 3.45
$ws.Range("P42").Value = 1.53
$ws.Range("Q42").Value = 2.1
$ws.Range("U42").Value = 9
$ws.Range("V42").Value = 8
$ws.Range("W42").Value = 22
$ws.Range("AA42").Value = 4.35
$ws.Range("AB42").Value = 12.5
$ws.Range("AE42").Value = 6.2
$ws.Range("AF42").Value = 14
$ws.Range("AI42").Value = 32
$ws.Range("G46").Value = 2.25
$ws.Range("I46").Value = 2.8
$ws.Range("L46").Value = 1.24
$ws.Range("U46").Value = 11.75
$ws.Range("W46").Value = 22
$ws.Range("X46").Value = 17.5
$ws.Range("AE46").Value = 10
$ws.Range("AF46").Value = 15
$ws.Range("AG46").Value = 10.25
$ws.Range("AH46").Value = 32
$ws.Range("AI46").Value = 22
$ws.Range("AJ46").Value = 29
$ws.Range("L47").Value = 1.23
$ws.Range("M47").Value = 3.4
$ws.Range("N47").Value = 1.7
$ws.Range("O47").Value = 1.93
$ws.Range("R47").Value = 1.57
$ws.Range("S47").Value = 2.1
$ws.Range("T47").Value = 8.5
$ws.Range("X47").Value = 14
$ws.Range("Z47").Value = 11.5
$ws.Range("AA47").Value = 6.7
$ws.Range("AB47").Value = 12.5
$ws.Range("AD47").Value = 300
$ws.Range("AE47").Value = 12.5
$ws.Range("AF47").Value = 22
$ws.Range("AJ47").Value = 32
$ws.Range("L48").Value = 1.17
$ws.Range("M48").Value = 5
$ws.Range("P59").Value = 1.4
$ws.Range("P60").Value = 1.4
$ws.Range("R63").Value = 1.67
$ws.Range("R64").Value = 1.53
$ws.Range("S64").Value = 2.38
$ws.Range("I65").Value = 3.1
$ws.Range("K65").Value = 12
$ws.Range("R65").Value = 1.62
$ws.Range("AE65").Value = 11
$ws.Range("AG65").Value = 11
$ws.Range("R66").Value = 1.7
$ws.Range("R67").Value = 1.83
$ws.Range("S67").Value = 1.83
$ws.Range("L69").Value = 1.25
$ws.Range("M69").Value = 3.75
$ws.Range("N69").Value = 1.9
$ws.Range("O69").Value = 1.95
$ws.Range("G76").Value = 1.52
$ws.Range("H76").Value = 3.95
$ws.Range("I76").Value = 5.7
$ws.Range("L76").Value = 1.27
$ws.Range("M76").Value = 3.1
$ws.Range("N76").Value = 1.82
$ws.Range("O76").Value = 1.8
$ws.Range("P76").Value = 1.38
$ws.Range("Q76").Value = 2.6
$ws.Range("R76").Value = 1.93
$ws.Range("S76").Value = 1.7
$ws.Range("T76").Value = 6.3
$ws.Range("U76").Value = 6.7
$ws.Range("W76").Value = 10.25
$ws.Range("X76").Value = 12.5
$ws.Range("Y76").Value = 30
$ws.Range("Z76").Value = 10.25
$ws.Range("AA76").Value = 7.8
$ws.Range("AB76").Value = 19.5
$ws.Range("AC76").Value = 100
$ws.Range("AD76").Value = 900
$ws.Range("AE76").Value = 14
$ws.Range("AF76").Value = 32
$ws.Range("AG76").Value = 18.5
$ws.Range("AH76").Value = 110
$ws.Range("AI76").Value = 65
$ws.Range("AJ76").Value = 65
$ws.Range("H78").Value = 3.3
$ws.Range("I78").Value = 4.05
$ws.Range("K78").Value = 6.8
$ws.Range("Q78").Value = 2.72
$ws.Range("T78").Value = 6.6
$ws.Range("X78").Value = 15
$ws.Range("Y78").Value = 28
$ws.Range("Z78").Value = 6.8
$ws.Range("AA78").Value = 6.5
$ws.Range("AE78").Value = 10.5
$ws.Range("AG78").Value = 14
$ws.Range("AJ78").Value = 50
$ws.Range("G80").Value = 2.9
$ws.Range("H80").Value = 3.05
$ws.Range("T80").Value = 8
$ws.Range("X80").Value = 28
$ws.Range("Y80").Value = 40
$ws.Range("AE80").Value = 7.2
$ws.Range("J84").Value = 1.06
$ws.Range("L84").Value = 1.36
$ws.Range("P84").Value = 1.44
$ws.Range("Q84").Value = 2.63
$ws.Range("R84").Value = 1.95
$ws.Range("S84").Value = 1.8
$ws.Range("J85").Value = 1.04
$ws.Range("L85").Value = 1.2
$ws.Range("P85").Value = 1.3
$ws.Range("R85").Value = 1.62
$ws.Range("J86").Value = 1.06
$ws.Range("L86").Value = 1.33
$ws.Range("P86").Value = 1.4
$ws.Range("R86").Value = 1.83
$ws.Range("S86").Value = 1.83
$ws.Range("J87").Value = 1.13
$ws.Range("L87").Value = 1.5
$ws.Range("P87").Value = 1.53
$ws.Range("Q87").Value = 2.38
$ws.Range("S87").Value = 1.67
